# Automatische test-sync: 2025-06-25 22:44:50
#
# Appends the next "Verzoek om factuur" test-mail log entry (row 6) to the
# "Logs" sheet, and bumps the matching "Factuur / Administratie" tally on
# the "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New log row, one past the current last row (row 5).
$newRow = 6

$logs.Cells.Item($newRow, 1).Value = "Verzoek om factuur"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Verzoek om factuur"
$logs.Cells.Item($newRow, 4).Value = "Factuur / Administratie"
$logs.Cells.Item($newRow, 5).Value = "nan"
$logs.Cells.Item($newRow, 6).Value = "2025-06-25 22:44:14"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# Dashboard!B2 tracks the "Factuur / Administratie" count; it was 3, now 4
# after the newly logged mail above.
$dashboard.Cells.Item(2, 2).Value = 4
